$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 510
$ws.Range("G3").Value = "不可售"
$ws.Range("F4").Value = 962
$ws.Range("F5").Value = 209
$ws.Range("F7").Value = 1064
$ws.Range("F8").Value = 847
$ws.Range("F9").Value = 253
$ws.Range("F11").Value = 81
$ws.Range("F12").Value = 846
$ws.Range("F13").Value = 293
$ws.Range("F14").Value = 588
$ws.Range("F15").Value = 507
$ws.Range("F16").Value = 1340
$ws.Range("F18").Value = 1267
$ws.Range("F19").Value = 1201
$ws.Range("F20").Value = 2894
$ws.Range("F21").Value = 1458
$ws.Range("F22").Value = 716
$ws.Range("F23").Value = 201
$ws.Range("F24").Value = 1282
$ws.Range("F26").Value = 1030
$ws.Range("F28").Value = 3160
$ws.Range("F29").Value = 615
$ws.Range("F31").Value = 1421
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 521
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 753
$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 753
$ws.Range("F4").Value = 510
$ws.Range("G4").Value = "不可售"
$ws.Range("F5").Value = 521
$ws.Range("F9").Value = 962
$ws.Range("F10").Value = 209
$ws.Range("F13").Value = 1064
$ws.Range("F14").Value = 847
$ws.Range("F15").Value = 253
$ws.Range("F21").Value = 81
$ws.Range("F23").Value = 846
$ws.Range("F24").Value = 293
$ws.Range("F25").Value = 588
$ws.Range("F26").Value = 507
$ws.Range("F27").Value = 1340
$ws.Range("F29").Value = 1267
$ws.Range("F30").Value = 1201
$ws.Range("F31").Value = 2894
$ws.Range("F32").Value = 1458
$ws.Range("F33").Value = 716
$ws.Range("F34").Value = 201
$ws.Range("F35").Value = 1282
$ws.Range("F39").Value = 1030
$ws.Range("F41").Value = 3160
$ws.Range("F42").Value = 615
$ws.Range("F44").Value = 1421